{"js": "// 1) Title: \"DOCUMENTO DE LAYOUT \u2013 CADASTRO DE USU\u00c1RIO (DOADOR)\"\n//           -> \"DOCUMENTO DE LAYOUT \u2013 CADASTRO DE USU\u00c1RIO (ONG)\"\n// Only \"(DOADOR)\" (with the parentheses) is targeted so the unrelated\n// \"DOADORES\" occurrence later in the document is left untouched.\nconst titleResults = context.document.body.search(\"(DOADOR)\", {\n  matchCase: true,\n  matchWholeWord: false\n});\ntitleResults.load(\"items\");\nawait context.sync();\n\nif (titleResults.items.length > 0) {\n  titleResults.items[0].insertText(\"(ONG)\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) Table cell: the date/time format explanation was previously split\n//    across several runs (with spell-check proofErr markers around\n//    \"dd\", \"yyyy\" and \"HH:mm:ss\"). Collapse it back into a single run\n//    that reads the same full sentence.\nconst dateTarget = 'Data e hora da gera\u00e7\u00e3o do arquivo, no formato \"dd-MM-yyyy HH:mm:ss\"';\nconst dateResults = context.document.body.search(dateTarget, {\n  matchCase: true\n});\ndateResults.load(\"items\");\nawait context.sync();\n\nif (dateResults.items.length > 0) {\n  const paragraph = dateResults.items[0].paragraphs.getFirst();\n  paragraph.insertText(dateTarget, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n$q = [char]34\n\n# --- Change 1: title \"... CADASTRO DE USUARIO (DOADOR)\" -> \"... (ONG)\" ---\n$rng1 = $d.Content\n$find1 = $rng1.Find\n$find1.ClearFormatting()\n$find1.MatchWildcards = $false\n$find1.Text = \"(DOADOR)\"\n$find1.MatchCase = $true\n$find1.MatchWholeWord = $false\n$found1 = $find1.Execute()\nif ($found1) {\n    $rng1.Text = \"(ONG)\"\n}\n\n# --- Change 2: merge the \"dd-MM-yyyy HH:mm:ss\" runs (plus spell-check markers)\n#     into a single run reading the full sentence, matching the cleaned-up cell. ---\n$rng2 = $d.Content\n$find2 = $rng2.Find\n$find2.ClearFormatting()\n$find2.MatchWildcards = $true\n$find2.Text = \"Data e hora da gera*o do arquivo, no formato \" + $q + \"dd-MM-yyyy HH:mm:ss\" + $q\n$find2.MatchCase = $true\n$found2 = $find2.Execute()\nif ($found2) {\n    # First collapse to a placeholder so the old multi-run/proofErr structure is\n    # discarded, then write the final sentence back as one clean run.\n    $rng2.Text = \"IRON_TMP_PLACEHOLDER\"\n    $rng2.Text = \"Data e hora da gera\u00e7\u00e3o do arquivo, no formato \" + $q + \"dd-MM-yyyy HH:mm:ss\" + $q\n}\n"}
